$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-12 Thursday", "2024-12-13 Friday"),
    @("30×84=2520", "63×86=5418"),
    @("16×38=608", "76×68=5168"),
    @("15×26=390", "17×28=476"),
    @("94×90=8460", "41×25=1025"),
    @("33×35=1155", "68×29=1972"),
    @("43×75=3225", "50×18=900"),
    @("43×52=2236", "68×28=1904"),
    @("55×72=3960", "81×16=1296"),
    @("34×66=2244", "23×97=2231"),
    @("31×40=1240", "46×62=2852"),
    @("81×96=7776", "68×96=6528"),
    @("63×57=3591", "82×98=8036"),
    @("23×36=828", "26×31=806"),
    @("51×92=4692", "26×81=2106"),
    @("35×14=490", "23×11=253"),
    @("45×13=585", "59×62=3658"),
    @("52×46=2392", "63×96=6048"),
    @("31×94=2914", "71×48=3408"),
    @("18×95=1710", "99×71=7029"),
    @("86×27=2322", "15×43=645"),
    @("15×72=1080", "81×54=4374"),
    @("77×75=5775", "48×62=2976"),
    @("92×77=7084", "80×13=1040"),
    @("89×93=8277", "47×36=1692"),
    @("70×28=1960", "31×65=2015")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
